$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the D-column formatting up front by copying each row's style
# from column C (D1/D3 get the shaded header style, D2/D4 get the plain
# bordered style), then write the actual values. ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# --- New column D: reducedTimePerMissionClear ---
$ws.Range("D1").Value = "reducedTimePerMissionClear"
$ws.Range("D2").Value = "(min)"
$ws.Range("D3").Value = "int"
$ws.Range("D4").Value = 3

# --- Existing cells that changed ---
# timeToExplosion's unit label switches from seconds to minutes ...
$ws.Range("C2").Value = "(min)"
# ... and its sample value drops from 20 to 15.
$ws.Range("C4").Value = 15

# --- Column D width ---
# NOTE: Excel's ColumnWidth (characters) differs from the stored OOXML
# <col width="..."> (MDW-rounded units) by a fixed offset for this sheet's
# font; 31 + 2/7 characters round-trips to a stored width of exactly 32.
$ws.Columns.Item(4).ColumnWidth = 31.285714285714285

# --- Selection moves as part of the edit session ---
$ws.Range("J11").Select()
